$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header shared-string text updates ---
$ws.Range("A8").Value = "Volume 31   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/29/2024  Through  5/5/2024"

# --- Number -> Text conversions (style 14, shared text "0" or "***.*") ---
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Formula = "=""0"""
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Formula = "=""0"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Formula = "=""***.*"""
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Formula = "=""0"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Formula = "=""***.*"""
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Formula = "=""0"""
$ws.Range("C27").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Formula = "=""0"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Formula = "=""***.*"""
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

# --- Text -> Number conversions ---
$ws.Range("D20").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 2
$ws.Range("H15").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("D20").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1
$ws.Range("D20").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 2
$ws.Range("H15").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = -50

# --- Plain value changes (same type) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -31.25
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 109
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = 17.204301075268
$ws.Range("L16").Value = 34.567901234567
$ws.Range("M16").Value = 9
$ws.Range("N16").Value = -83.204930662557
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -18.181818181818
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -12.820512820512
$ws.Range("I17").Value = 191
$ws.Range("J17").Value = 158
$ws.Range("K17").Value = 20.886075949367
$ws.Range("L17").Value = 7.909604519774
$ws.Range("M17").Value = 122.093023255814
$ws.Range("N17").Value = -26.819923371647
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = -2.597402597402
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -56.140350877193
$ws.Range("N18").Value = -89.035087719298
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 80
$ws.Range("H19").Value = -15
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 329
$ws.Range("K19").Value = -1.215805471124
$ws.Range("L19").Value = 14.035087719298
$ws.Range("M19").Value = 43.171806167400
$ws.Range("N19").Value = -14.021164021164
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 73.684210526315
$ws.Range("I20").Value = 123
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = 12.844036697247
$ws.Range("L20").Value = 53.75
$ws.Range("N20").Value = -87.332646755921
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = 1.960784313725
$ws.Range("F21").Value = 177
$ws.Range("G21").Value = 184
$ws.Range("H21").Value = -3.804347826086
$ws.Range("I21").Value = 834
$ws.Range("J21").Value = 781
$ws.Range("K21").Value = 6.786171574903
$ws.Range("L21").Value = 14.718019257221
$ws.Range("M21").Value = 23.190546528803
$ws.Range("N21").Value = -71.881321645313
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -11.111111111111
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -32
$ws.Range("I23").Value = 67
$ws.Range("J23").Value = 79
$ws.Range("K23").Value = -15.189873417721
$ws.Range("L23").Value = -14.102564102564
$ws.Range("M23").Value = 48.888888888888
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = 24.390243902439
$ws.Range("F24").Value = 156
$ws.Range("H24").Value = -8.771929824561
$ws.Range("I24").Value = 774
$ws.Range("J24").Value = 795
$ws.Range("K24").Value = -2.641509433962
$ws.Range("L24").Value = -10.416666666666
$ws.Range("M24").Value = 54.183266932270
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 121.428571428571
$ws.Range("F25").Value = 85
$ws.Range("G25").Value = 83
$ws.Range("H25").Value = 2.409638554216
$ws.Range("I25").Value = 417
$ws.Range("J25").Value = 388
$ws.Range("K25").Value = 7.474226804123
$ws.Range("L25").Value = 31.545741324921
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 46.666666666666
$ws.Range("F26").Value = 72
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = 28.571428571428
$ws.Range("I26").Value = 311
$ws.Range("J26").Value = 263
$ws.Range("K26").Value = 18.250950570342
$ws.Range("L26").Value = -3.715170278637
$ws.Range("M26").Value = 4.713804713804
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("L27").Value = 12.5
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -15.625
$ws.Range("L28").Value = 22.727272727272
$ws.Range("N29").Value = -76.923076923076
$ws.Range("N30").Value = -84.615384615384
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = -33.333333333333

$excel.CutCopyMode = $false
